$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Note" column (D) for the 3D-printed parts (rows 2-9) is being
# reworded from the generic "3DP" to the more specific "3DP(PET)".
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = "3DP(PET)"
}

# Column D's best-fit width shrinks now that the longer header text no
# longer dominates the autofit calculation.
$ws.Columns.Item(4).ColumnWidth = 9.71

# Selection moves to cover the populated BOM table.
$ws.Range("A1:D16").Select() | Out-Null
